$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.00006486019690155054
$ws.Range("C2").Value = 0.05231270169004087
$ws.Range("D2").Value = 0.7127328510149897
$ws.Range("E2").Value = 71517.89157740913
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 71518.65668782203
